$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2021" year column (R) mirroring the formatting of the
# existing "2020" year column (Q4 header cell).
$null = $ws.Range("Q4").Copy()
$null = $ws.Range("R4").PasteSpecial(-4122)
$ws.Range("R4").Value = 2021

# Add the corresponding data value for 2021 in R5, using the same
# number formatting (0.0) as other data cells such as H5/M5.
$null = $ws.Range("H5").Copy()
$null = $ws.Range("R5").PasteSpecial(-4122)
$ws.Range("R5").Value = 18.953297329007047

# Update the active selection to reflect the new cursor location.
$null = $ws.Range("Q8").Select()
